$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original text representation (no numeric
# auto-conversion / scientific notation / trailing-zero loss) by forcing the
# cell format to Text before writing the new value.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.892.23"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.809.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4285"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3684"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07241"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8631"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "2.047.36"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +17.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.15"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.620"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.391"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06926"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008840"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.24"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.948.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.193"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.93"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.280.55"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +15.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.09"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.885"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.37"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.229"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.908"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +15.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.70"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08941"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7393"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.49%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.423"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.807"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.122"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05218"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01921"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5079"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.755"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +12.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1648"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.447"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.306"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "107.19"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.38"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.64%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4578"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.96%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.647"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06272"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.810"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.13%  "
